# Scheduled market-price refresh: updates cached Universalis price/profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1949.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1949.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1949.5
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -2299.5
$ws.Range("H40").Value = 6741.6665
$ws.Range("I40").Value = 5180
$ws.Range("J40").Value = 7857.143
$ws.Range("K40").Value = 5180
$ws.Range("L40").Value = 7857.143
$ws.Range("M40").Value = -5005
$ws.Range("N40").Value = -8207.143
$ws.Range("H46").Value = 2208.5
$ws.Range("I46").Value = 2208.5
$ws.Range("K46").Value = 6625.5
$ws.Range("M46").Value = -6506.5
$ws.Range("H48").Value = 1019.8
$ws.Range("I48").Value = 875
$ws.Range("J48").Value = 1599
$ws.Range("K48").Value = 2625
$ws.Range("L48").Value = 4797
$ws.Range("M48").Value = -2333
$ws.Range("N48").Value = -5381
$ws.Range("H56").Value = 1019.8
$ws.Range("I56").Value = 875
$ws.Range("J56").Value = 1599
$ws.Range("K56").Value = 2625
$ws.Range("L56").Value = 4797
$ws.Range("M56").Value = -2091
$ws.Range("N56").Value = -5865
$ws.Range("H60").Value = 2208.5
$ws.Range("I60").Value = 2208.5
$ws.Range("K60").Value = 6625.5
$ws.Range("M60").Value = -6141.5
$ws.Range("H62").Value = 2226
$ws.Range("I62").Value = 1829.7142
$ws.Range("K62").Value = 1829.7142
$ws.Range("M62").Value = -1205.7142
$ws.Range("H65").Value = 2226
$ws.Range("I65").Value = 1829.7142
$ws.Range("K65").Value = 9148.571
$ws.Range("M65").Value = -6028.571
$ws.Range("H113").Value = 9271.933999999999
$ws.Range("J113").Value = 10177.667
$ws.Range("L113").Value = 10177.667
$ws.Range("N113").Value = -16685.667
$ws.Range("H132").Value = 25170.2
$ws.Range("I132").Value = 3963.0588
$ws.Range("J132").Value = 70235.375
$ws.Range("K132").Value = 11889.1764
$ws.Range("L132").Value = 210706.125
$ws.Range("M132").Value = -9359.1764
$ws.Range("N132").Value = -215766.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 53358.168
$ws.Range("J2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("N2").Value = -4226
$ws.Range("H61").Value = 1619.45
$ws.Range("I61").Value = 1547.3158
$ws.Range("J61").Value = 2990
$ws.Range("K61").Value = 1547.3158
$ws.Range("L61").Value = 2990
$ws.Range("M61").Value = -1335.3158
$ws.Range("N61").Value = -3414
$ws.Range("H74").Value = 2413.56
$ws.Range("I74").Value = 2188.652
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2188.652
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1314.652
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 2413.56
$ws.Range("I77").Value = 2188.652
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 10943.26
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -6575.26
$ws.Range("N77").Value = -33736
$ws.Range("H96").Value = 64672
$ws.Range("J96").Value = 64672
$ws.Range("L96").Value = 64672
$ws.Range("N96").Value = -70164
$ws.Range("H110").Value = 6502.25
$ws.Range("I110").Value = 6502.25
$ws.Range("K110").Value = 6502.25
$ws.Range("M110").Value = -4457.25
$ws.Range("H116").Value = 53358.168
$ws.Range("J116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("N116").Value = -8588
$ws.Range("H122").Value = 2023.0555
$ws.Range("I122").Value = 2023.0555
$ws.Range("K122").Value = 6069.166499999999
$ws.Range("M122").Value = -3619.166499999999
$ws.Range("H132").Value = 2407.7856
$ws.Range("I132").Value = 2439.1538
$ws.Range("K132").Value = 7317.4614
$ws.Range("M132").Value = -4787.4614
$ws.Range("H136").Value = 1619.45
$ws.Range("I136").Value = 1547.3158
$ws.Range("J136").Value = 2990
$ws.Range("K136").Value = 4641.9474
$ws.Range("L136").Value = 8970
$ws.Range("M136").Value = -2091.9474
$ws.Range("N136").Value = -14070

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 53358.168
$ws.Range("J3").Value = 4000
$ws.Range("L3").Value = 4000
$ws.Range("N3").Value = -4228
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H86").Value = 2335.875
$ws.Range("I86").Value = 1925.4375
$ws.Range("J86").Value = 3156.75
$ws.Range("K86").Value = 1925.4375
$ws.Range("L86").Value = 3156.75
$ws.Range("M86").Value = -802.4375
$ws.Range("N86").Value = -5402.75
$ws.Range("H88").Value = 44892
$ws.Range("J88").Value = 44892
$ws.Range("L88").Value = 44892
$ws.Range("N88").Value = -45704
$ws.Range("H89").Value = 2335.875
$ws.Range("I89").Value = 1925.4375
$ws.Range("J89").Value = 3156.75
$ws.Range("K89").Value = 9627.1875
$ws.Range("L89").Value = 15783.75
$ws.Range("M89").Value = -4011.1875
$ws.Range("N89").Value = -27015.75
$ws.Range("H91").Value = 44892
$ws.Range("J91").Value = 44892
$ws.Range("L91").Value = 44892
$ws.Range("N91").Value = -47700
$ws.Range("H134").Value = 1318.1082
$ws.Range("I134").Value = 1318.1082
$ws.Range("K134").Value = 3954.3246
$ws.Range("M134").Value = -1419.3246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1652.2
$ws.Range("I31").Value = 1533.7
$ws.Range("J31").Value = 2126.2
$ws.Range("K31").Value = 1533.7
$ws.Range("L31").Value = 2126.2
$ws.Range("M31").Value = -1238.7
$ws.Range("N31").Value = -2716.2
$ws.Range("H34").Value = 1652.2
$ws.Range("I34").Value = 1533.7
$ws.Range("J34").Value = 2126.2
$ws.Range("K34").Value = 1533.7
$ws.Range("L34").Value = 2126.2
$ws.Range("M34").Value = -1331.7
$ws.Range("N34").Value = -2530.2
$ws.Range("H58").Value = 3366.7693
$ws.Range("I58").Value = 2785.3333
$ws.Range("J58").Value = 4675
$ws.Range("K58").Value = 2785.3333
$ws.Range("L58").Value = 4675
$ws.Range("M58").Value = -2582.3333
$ws.Range("N58").Value = -5081
$ws.Range("H94").Value = 10624.818
$ws.Range("J94").Value = 2271
$ws.Range("L94").Value = 2271
$ws.Range("N94").Value = -3173
$ws.Range("H132").Value = 1948.8096
$ws.Range("I132").Value = 1878.8948
$ws.Range("J132").Value = 2613
$ws.Range("K132").Value = 5636.6844
$ws.Range("L132").Value = 7839
$ws.Range("M132").Value = -3106.6844
$ws.Range("N132").Value = -12899
$ws.Range("H134").Value = 2646.5
$ws.Range("I134").Value = 2118.074
$ws.Range("K134").Value = 6354.222
$ws.Range("M134").Value = -3819.222
$ws.Range("H136").Value = 3366.7693
$ws.Range("I136").Value = 2785.3333
$ws.Range("J136").Value = 4675
$ws.Range("K136").Value = 8355.999899999999
$ws.Range("L136").Value = 14025
$ws.Range("M136").Value = -5805.999899999999
$ws.Range("N136").Value = -19125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2587.5
$ws.Range("I60").Value = 1000
$ws.Range("J60").Value = 3116.6667
$ws.Range("K60").Value = 3000
$ws.Range("L60").Value = 9350.000100000001
$ws.Range("M60").Value = -2749
$ws.Range("N60").Value = -9852.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 45141.465
$ws.Range("I80").Value = 78085.336
$ws.Range("K80").Value = 78085.336
$ws.Range("M80").Value = -77087.336
$ws.Range("H83").Value = 45141.465
$ws.Range("I83").Value = 78085.336
$ws.Range("K83").Value = 390426.68
$ws.Range("M83").Value = -385434.68
$ws.Range("H96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -50492
$ws.Range("H102").Value = 3659.16
$ws.Range("I102").Value = 2432.7334
$ws.Range("J102").Value = 5498.8
$ws.Range("K102").Value = 2432.7334
$ws.Range("L102").Value = 5498.8
$ws.Range("M102").Value = -810.7334000000001
$ws.Range("N102").Value = -8742.799999999999
$ws.Range("H132").Value = 3805.818
$ws.Range("I132").Value = 3758
$ws.Range("J132").Value = 3933.3333
$ws.Range("K132").Value = 11274
$ws.Range("L132").Value = 11799.9999
$ws.Range("M132").Value = -8744
$ws.Range("N132").Value = -16859.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5954.423
$ws.Range("I7").Value = 5902.353
$ws.Range("K7").Value = 5902.353
$ws.Range("M7").Value = -5790.353
$ws.Range("H93").Value = 7939.4736
$ws.Range("I93").Value = 7740.5835
$ws.Range("J93").Value = 8280.429
$ws.Range("K93").Value = 7740.5835
$ws.Range("L93").Value = 8280.429
$ws.Range("M93").Value = -6492.5835
$ws.Range("N93").Value = -10776.429
$ws.Range("H126").Value = 5954.423
$ws.Range("I126").Value = 5902.353
$ws.Range("K126").Value = 17707.059
$ws.Range("M126").Value = -15237.059
$ws.Range("H132").Value = 2995.5667
$ws.Range("I132").Value = 2180.4092
$ws.Range("K132").Value = 6541.2276
$ws.Range("M132").Value = -4011.2276
$ws.Range("H136").Value = 2517.2334
$ws.Range("I136").Value = 2580.6206
$ws.Range("K136").Value = 7741.861800000001
$ws.Range("M136").Value = -5191.861800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2066.5881
$ws.Range("I122").Value = 1918.6923
$ws.Range("J122").Value = 2547.25
$ws.Range("K122").Value = 5756.0769
$ws.Range("L122").Value = 7641.75
$ws.Range("M122").Value = -3306.0769
$ws.Range("N122").Value = -12541.75
$ws.Range("H132").Value = 2008.7894
$ws.Range("I132").Value = 1916.875
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 5750.625
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -3220.625
$ws.Range("N132").Value = -12557
$ws.Range("H136").Value = 1552.4286
$ws.Range("I136").Value = 890.43475
$ws.Range("K136").Value = 2671.30425
$ws.Range("M136").Value = -121.3042500000001
